{"js": "// Commit message: \"Added PDF versions to site\"\n//\n// The \"Additional resources\" bullet list used to start with a bullet\n// linking to a Word version of this document (a hyperlink whose target was\n// a /word-versions/*.docx file). That bullet is being dropped entirely \u2014\n// the list now starts directly with \"Instructor orientation\". Everything\n// else in the document (headings, bookmarks, the remaining hyperlinks)\n// stays exactly the same; only this one paragraph disappears.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst TARGET_TEXT = \"Word version of this document\";\n\n// Walk back-to-front so deleting a paragraph never invalidates the index of\n// one we haven't visited yet.\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === TARGET_TEXT) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Commit message: \"Added PDF versions to site\"\n#\n# The \"Additional resources\" bullet list used to start with a bullet linking\n# to a Word version of this document (a hyperlink whose target was a\n# /word-versions/*.docx file). That bullet is being dropped entirely -- the\n# list now starts directly with \"Instructor orientation\". Everything else in\n# the document (headings, bookmarks, the remaining hyperlinks) stays exactly\n# the same; only this one paragraph disappears.\n\n$d = $word.ActiveDocument\n$target = \"Word version of this document\"\n\n# Walk back-to-front so deleting a paragraph never invalidates the 1-based\n# index of one we haven't visited yet.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n    if ($text -eq $target) {\n        $p.Range.Delete()\n    }\n}\n"}
